$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feature")

# ---------------------------------------------------------------------------
# 0) Capture the two pre-existing cell styles we'll need to reuse later,
#    into scratch cells far off to the right, *before* anything else in this
#    script mutates their source cells (E3/E4 currently carry the "open"
#    themed-fill style; F19 carries the date-format style).
# ---------------------------------------------------------------------------
$ws.Range("E3").Copy($ws.Range("Z1"))     # scratch: themed-fill "open" style (cellXf 3)
$ws.Range("Z1").ClearContents()
$ws.Range("F19").Copy($ws.Range("Z2"))    # scratch: date-format style (cellXf 2)
$ws.Range("Z2").ClearContents()

# ---------------------------------------------------------------------------
# 1) Structural: insert the three new rows that the diff introduces.
#    Doing this top-down (inserting lower rows after higher ones have already
#    pushed everything down) keeps the target row numbers stable.
# ---------------------------------------------------------------------------
$ws.Rows("6:6").Insert()    # new "view search results by different criteria" row
$ws.Rows("20:20").Insert()  # blank spacer row before "insert book"
$ws.Rows("24:24").Insert()  # new "view user" row

# ---------------------------------------------------------------------------
# 2) Row 3 ("add photo"): status flips from "open"(themed fill) to "done",
#    gains a date and a remark.
# ---------------------------------------------------------------------------
$ws.Range("E3").ClearFormats()
$ws.Range("E3").Value = "done"
$ws.Range("Z2").Copy($ws.Range("F3"))
$ws.Range("F3").Value = 43983
$ws.Range("H3").Value = "/ebook/images/"

# ---------------------------------------------------------------------------
# 3) Row 4 ("show data in multipages"): stays "open" but now uses the new
#    orange highlight fill (fillId 4); gains a dated cell using the same new
#    fill plus the date number format (fillId 4 + numFmt 16).
# ---------------------------------------------------------------------------
$ws.Range("E4").Interior.Color = 49407
$ws.Range("E4").Value = "open"
$ws.Range("F4").Interior.Color = 49407
$ws.Range("F4").NumberFormat = "d-mmm"
$ws.Range("F4").Value = 43984

# ---------------------------------------------------------------------------
# 4) New row 6: "view search results by different criteria" / open, reusing
#    the pre-existing themed "open" fill style captured in Z1.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy($ws.Range("E6"))
$ws.Range("E6").Value = "open"
$ws.Range("C6").Value = "view search results by different criteria"

# ---------------------------------------------------------------------------
# 5) Rows 16/17 (order / view order): the "issue" note moves up from the
#    "view order" row to the "order" row, leaving the old cell's highlighted
#    style behind but empty.
# ---------------------------------------------------------------------------
$ws.Range("E17").Copy($ws.Range("E16"))
$ws.Range("E16").Value = "issue ,can not populate order data into db"
$ws.Range("E17").ClearContents()

# ---------------------------------------------------------------------------
# 6) Row 20 (blank spacer): keep the date-format style, value empty.
# ---------------------------------------------------------------------------
$ws.Range("Z2").Copy($ws.Range("F20"))

# ---------------------------------------------------------------------------
# 7) Rows 21-23 (insert book / update book / delete book): add "done"/"open"
#    status + dates that didn't exist before.
# ---------------------------------------------------------------------------
$ws.Range("E21").Value = "done"
$ws.Range("E22").Value = "done"
$ws.Range("Z2").Copy($ws.Range("F22"))
$ws.Range("F22").Value = 43983

$ws.Range("Z1").Copy($ws.Range("E23"))
$ws.Range("E23").Value = "open"

# ---------------------------------------------------------------------------
# 8) New row 24: "view user" with an empty themed-style status cell and an
#    empty date-style cell (mirrors row 23/25's blank date cells).
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = "view user"
$ws.Range("Z1").Copy($ws.Range("E24"))
$ws.Range("Z2").Copy($ws.Range("F24"))

# ---------------------------------------------------------------------------
# 9) Row 23 and 25 also need the plain date-format (empty) cell in F.
# ---------------------------------------------------------------------------
$ws.Range("Z2").Copy($ws.Range("F23"))
$ws.Range("Z2").Copy($ws.Range("F25"))

# ---------------------------------------------------------------------------
# 10) Rows 25-27 (search book by criteira / search order / process invoice):
#     add "open" status cells using the same themed fill.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy($ws.Range("E25"))
$ws.Range("E25").Value = "open"
$ws.Range("Z1").Copy($ws.Range("E26"))
$ws.Range("E26").Value = "open"
$ws.Range("Z1").Copy($ws.Range("E27"))
$ws.Range("E27").Value = "open"

# ---------------------------------------------------------------------------
# 11) Brand-new rows 28-33: email / post / review features, each with a
#     label row and a blank "open" status spacer row underneath, matching
#     the pattern already used above.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy($ws.Range("E28"))
$ws.Range("E28").Value = "open"

$ws.Range("B29").Value = "email"
$ws.Range("Z1").Copy($ws.Range("E29"))
$ws.Range("E29").Value = "open"

$ws.Range("Z1").Copy($ws.Range("E30"))
$ws.Range("E30").Value = "open"

$ws.Range("B31").Value = "post"
$ws.Range("Z1").Copy($ws.Range("E31"))
$ws.Range("E31").Value = "open"

$ws.Range("Z1").Copy($ws.Range("E32"))
$ws.Range("E32").Value = "open"

$ws.Range("B33").Value = "review"
$ws.Range("Z1").Copy($ws.Range("E33"))
$ws.Range("E33").Value = "open"

# ---------------------------------------------------------------------------
# 12) Clean up the scratch cells used as style donors.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Clear()
$ws.Range("Z2").Clear()

# ---------------------------------------------------------------------------
# 13) Column widths / layout tweaks from the diff.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 3.1640625
$ws.Columns("C:C").ColumnWidth = 26.33203125
$ws.Columns("D:D").ColumnWidth = 10.5
$ws.Columns("E:E").ColumnWidth = 8.6640625

# Selection bookkeeping to mirror the saved workbook state.
$ws.Range("C22").Select()
